$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.314.57"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "3.153.47"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'591.60"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'147.61"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.152.75"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("E10").Value = "  -2.99%  "
$ws.Range("D11").Value = "'5.89"
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("D14").Value = "'37.24"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("D15").Value = "3.675.56"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.291.17"
$ws.Range("E17").Value = "  +5.18%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.974.58"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "'7.20"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'467.78"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'14.38"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").Value = "'0.735"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'7.45"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").Value = "'2.36"
$ws.Range("E24").Value = "  +9.08%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'13.05"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "'81.10"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "'9.69"
$ws.Range("E28").Value = "  +11.77%  "
$ws.Range("D29").Value = "'2.71"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "'7.25"
$ws.Range("E32").Value = "  +6.17%  "
$ws.Range("D33").Value = "'27.69"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D35").Value = "0.0₃0850"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("D37").Value = "'2.34"
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("D38").Value = "'6.09"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("E39").Value = "  -4.12%  "
$ws.Range("D40").Value = "'51.66"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "'456.46"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "'9.09"
$ws.Range("E42").Value = "  +4.01%  "
$ws.Range("D43").Value = "'0.293"
$ws.Range("E43").Value = "  +5.57%  "
$ws.Range("D44").Value = "'0.0373"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "2.921.03"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "'39.59"
$ws.Range("E46").Value = "  +14.10%  "
$ws.Range("D47").Value = "'0.109"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").Value = "'127.16"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D50").Value = "'2.25"
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("E51").Value = "  -0.56%  "
